$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45175
}
